$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matching source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values per row
$ws.Range("D2").Value = '30.382.23'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.875.39'
$ws.Range("E3").Value = '  -0.87%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '238.44'
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '0.4806'
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("E8").Value = '  -2.96%  '
$ws.Range("D9").Value = '0.06504'
$ws.Range("E9").Value = '  -1.54%  '
$ws.Range("D10").Value = '1.871.07'
$ws.Range("E10").Value = '  -1.22%  '
$ws.Range("D11").Value = '0.07485'
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("E12").Value = '  -2.20%  '
$ws.Range("D13").Value = '5.065'
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("D14").Value = '88.23'
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("D15").Value = '0.6600'
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("D16").Value = '30.354.12'
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").Value = '13.27'
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = '0.000007572'
$ws.Range("E19").Value = '  -2.58%  '
$ws.Range("D20").Value = '2.113.72'
$ws.Range("E20").Value = '  -1.25%  '
$ws.Range("D21").Value = '5.290'
$ws.Range("E21").Value = '  -3.80%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").Value = '220.16'
$ws.Range("E23").Value = '  +14.60%  '
$ws.Range("D24").Value = '6.173'
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("D25").Value = '9.327'
$ws.Range("E25").Value = '  -1.38%  '
$ws.Range("D26").Value = '166.40'
$ws.Range("E26").Value = '  +0.79%  '
$ws.Range("D27").Value = '18.42'
$ws.Range("D28").Value = '1.960'
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("D29").Value = '1.462'
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("D30").Value = '0.09365'
$ws.Range("E30").Value = '  +2.26%  '
$ws.Range("D31").Value = '4.298'
$ws.Range("E31").Value = '  +0.79%  '
$ws.Range("D32").Value = '4.018'
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("D33").Value = '0.05017'
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("D34").Value = '1.201'
$ws.Range("E34").Value = '  +4.98%  '
$ws.Range("D35").Value = '0.7418'
$ws.Range("E35").Value = '  +0.90%  '
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("D38").Value = '2.613'
$ws.Range("E38").Value = '  -1.31%  '
$ws.Range("D39").Value = '0.9052'
$ws.Range("E39").Value = '  -1.67%  '
$ws.Range("D40").Value = '2.058'
$ws.Range("E40").Value = '  -1.28%  '
$ws.Range("D41").Value = '106.68'
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").Value = '5.860'
$ws.Range("E42").Value = '  -0.87%  '
$ws.Range("D43").Value = '0.4261'
$ws.Range("E43").Value = '  -1.93%  '
$ws.Range("E44").Value = '  +0.23%  '
$ws.Range("E45").Value = '  -3.14%  '
$ws.Range("D46").Value = '64.09'
$ws.Range("E46").Value = '  -2.17%  '
$ws.Range("D47").Value = '0.1270'
$ws.Range("E47").Value = '  -7.68%  '
$ws.Range("D48").Value = '1.470'
$ws.Range("E48").Value = '  -6.46%  '
$ws.Range("D49").Value = '8.895'
$ws.Range("E49").Value = '  -1.53%  '
$ws.Range("E50").Value = '  -1.95%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05631'
$ws.Range("E51").Value = '  -2.61%  '
